$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 255
$ws.Range("F3").Value = 303
$ws.Range("F4").Value = 2941
$ws.Range("F5").Value = 74
$ws.Range("F8").Value = 1610
$ws.Range("F10").Value = 842
$ws.Range("F12").Value = 2638
$ws.Range("F14").Value = 1482
$ws.Range("F15").Value = 6935
$ws.Range("F17").Value = 7095
$ws.Range("F19").Value = 2965
$ws.Range("F21").Value = 3447
$ws.Range("F23").Value = 149
$ws.Range("F24").Value = 1822
$ws.Range("F26").Value = 290
$ws.Range("F27").Value = 869
$ws.Range("F28").Value = 10
$ws.Range("F33").Value = 2536
$ws.Range("F35").Value = 159
$ws.Range("F36").Value = 367
$ws.Range("F37").Value = 1010
$ws.Range("F38").Value = 204
$ws.Range("F39").Value = 457
$ws.Range("F40").Value = 508

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 6
$ws.Range("F13").Value = 47
$ws.Range("F18").Value = 4

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 303
$ws.Range("F6").Value = 2941
$ws.Range("F7").Value = 74
$ws.Range("F9").Value = 1610
$ws.Range("F11").Value = 842
$ws.Range("F14").Value = 2638
$ws.Range("F15").Value = 1482
$ws.Range("F18").Value = 6
$ws.Range("F20").Value = 6935
$ws.Range("F22").Value = 7095
$ws.Range("F24").Value = 2966
$ws.Range("F26").Value = 3447
$ws.Range("F29").Value = 47
$ws.Range("F31").Value = 1822
$ws.Range("F34").Value = 290
$ws.Range("F35").Value = 869
$ws.Range("F36").Value = 10
$ws.Range("F41").Value = 2536
$ws.Range("F43").Value = 159
$ws.Range("F44").Value = 4
$ws.Range("F45").Value = 367
$ws.Range("F46").Value = 1010
$ws.Range("F47").Value = 204
$ws.Range("F48").Value = 457
$ws.Range("F49").Value = 508
